$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9 (shifts existing rows 9-48 down to 10-49),
# inheriting the formatting of the row above (row 8).
$ws.Rows.Item(9).Insert(-4121)

# Fill in the new "Event7 / add 1005 new character" row.
$ws.Range("A9").Value = 107
$ws.Range("B9").Value = "Event7"
$ws.Range("C9").Value = "add 1005 new character"
$ws.Range("D9").Value = 23
$ws.Range("E9").Value = "[108]"
$ws.Range("F9").Value = "Tutorial_01"

# Move the selection to F9 (was F10 before the insert).
$ws.Range("F9").Select() | Out-Null
